# Update cryptocurrency price/volume table (columns D and E) for rows 2-51
# to match the refreshed data feed, per the GitHub Actions scheduled update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.786.26"
$ws.Range("D2").Style = $style
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.408.16"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.04%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.01"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.68%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.85"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E8").Value = "  +3.79%  "
$ws.Range("E9").Value = "  -1.93%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.68"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("E11").Value = "  -1.13%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  +2.42%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.835.52"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.57%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.727.72"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -2.08%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.423.68"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.16%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.31"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.41%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.41"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.85%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.38"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -1.89%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("E22").Value = "  +0.02%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.48"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.11%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.37"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -1.00%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0770"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -2.47%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.49"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.08%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("E38").Value = "  -2.00%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "319.06"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +2.03%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.405"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -4.24%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.50%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.43"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -2.19%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0967"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.09%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.53"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("E45").Value = "  -1.85%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.576"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.76%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0223"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.57%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.385"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("E49").Value = "  -2.74%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.03"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -3.52%  "
